# Atualiza os valores da planilha "Valores" (relatorio_neomater_APENAS_VALORES)
# conforme nova extracao de dados - cria backup do arquivo que alimenta o PBI.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores")

# Linha 2
$ws.Range("A2").Value = 4
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 14
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 6

# Linha 9
$ws.Range("A9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = 5

# Linha 10
$ws.Range("A10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = 9
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2

# Linha 11
$ws.Range("A11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("H11").Value = 0

# Linha 13
$ws.Range("A13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("L13").Value = 0

# Linha 15
$ws.Range("A15").Value = 3
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 3
$ws.Range("J15").Value = 4
$ws.Range("L15").Value = 5

# Linha 16 (totais)
$ws.Range("A16").Value = 11
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 31
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 13
